$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New facility row 5: "Crumpet GEF" ---------------------------------
# Copy the formatting from the row above (row 4, the last populated data
# row) down into the previously-unused row 5 so the new row matches the
# look of the rest of the table.
$ws.Range("A4:J4").Copy() | Out-Null
$ws.Range("A5:J5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Rows.Item(5).RowHeight = $ws.Rows.Item(4).RowHeight

$ws.Range("A5").Value = "Crumpet GEF"
$ws.Range("B5").Value = 20001371
$ws.Range("C5").Value = "Crumpet exporter"
$ws.Range("D5").Value = "GBP"
$ws.Range("E5").Value = 7000000
$ws.Range("F5").Value = 3938753.8
$ws.Range("G5").Value = 777
$ws.Range("H5").Value = 456
$ws.Range("I5").Value = "GBP"
$ws.Range("J5").Value = "GBP"

# --- Row 6 (previously blank placeholder): "Scone GEF" ------------------
$ws.Range("A6").Value = "Scone GEF"
$ws.Range("B6").Value = 20001371
$ws.Range("C6").Value = "Scone exporter"
$ws.Range("D6").Value = "GBP"
$ws.Range("E6").Value = 770000
$ws.Range("F6").Value = 761579.37
$ws.Range("G6").Value = 777
$ws.Range("H6").Value = 456.77
$ws.Range("I6").Value = "GBP"
$ws.Range("J6").Value = "GBP"

# --- Restore the view: scroll back to column A and select D7 ------------
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("D7").Select() | Out-Null

$wb.Save()
